$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all existing data rows (2-385)
# from serial 45192 (2023-09-23) to serial 45202 (2023-10-03)
$ws.Range("C2:C385").Value = 45202

# Row 385 gains an explicit row height (ht="15" customHeight="1") in the new file
$ws.Rows.Item(385).RowHeight = 15

# --- Add new row 386 ---
$ws.Range("A386").Value = "A 45395-2023"
$ws.Range("B386").Value = 45194
$ws.Range("B386").NumberFormat = "YYYY-MM-DD"
$ws.Range("C386").Value = 45202
$ws.Range("C386").NumberFormat = "YYYY-MM-DD"
$ws.Range("D386").Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Range("E386").Value = "SVENLJUNGA"
$ws.Range("G386").Value = 2.4
$ws.Range("H386").Value = 0
$ws.Range("I386").Value = 0
$ws.Range("J386").Value = 0
$ws.Range("K386").Value = 0
$ws.Range("L386").Value = 0
$ws.Range("M386").Value = 0
$ws.Range("N386").Value = 0
$ws.Range("O386").Value = 0
$ws.Range("P386").Value = 0
$ws.Range("Q386").Value = 0
$ws.Range("R386").Value = ""
$ws.Range("R386").WrapText = $true
$ws.Rows.Item(386).RowHeight = 15

# --- Add new row 387 ---
$ws.Range("A387").Value = "A 45400-2023"
$ws.Range("B387").Value = 45194
$ws.Range("B387").NumberFormat = "YYYY-MM-DD"
$ws.Range("C387").Value = 45202
$ws.Range("C387").NumberFormat = "YYYY-MM-DD"
$ws.Range("D387").Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Range("E387").Value = "SVENLJUNGA"
$ws.Range("G387").Value = 1.7
$ws.Range("H387").Value = 0
$ws.Range("I387").Value = 0
$ws.Range("J387").Value = 0
$ws.Range("K387").Value = 0
$ws.Range("L387").Value = 0
$ws.Range("M387").Value = 0
$ws.Range("N387").Value = 0
$ws.Range("O387").Value = 0
$ws.Range("P387").Value = 0
$ws.Range("Q387").Value = 0
$ws.Range("R387").Value = ""
$ws.Range("R387").WrapText = $true
$ws.Rows.Item(387).RowHeight = 15

# --- Add new row 388 (no explicit row height set, matching the source diff) ---
$ws.Range("A388").Value = "A 45397-2023"
$ws.Range("B388").Value = 45194
$ws.Range("B388").NumberFormat = "YYYY-MM-DD"
$ws.Range("C388").Value = 45202
$ws.Range("C388").NumberFormat = "YYYY-MM-DD"
$ws.Range("D388").Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Range("E388").Value = "SVENLJUNGA"
$ws.Range("G388").Value = 1.8
$ws.Range("H388").Value = 0
$ws.Range("I388").Value = 0
$ws.Range("J388").Value = 0
$ws.Range("K388").Value = 0
$ws.Range("L388").Value = 0
$ws.Range("M388").Value = 0
$ws.Range("N388").Value = 0
$ws.Range("O388").Value = 0
$ws.Range("P388").Value = 0
$ws.Range("Q388").Value = 0
$ws.Range("R388").Value = ""
$ws.Range("R388").WrapText = $true
